$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample import data (rows 2-8, columns A-B) but keep cell formatting.
$ws.Range("A2:B8").ClearContents()

# Update the active selection to reflect the cleared state.
$ws.Range("E14").Select()
